$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Q2 and Q3 values from 5 to 0
$ws.Range("Q2").Value = 0
$ws.Range("Q3").Value = 0

# Update the selection to Q2:Q3 with active cell Q2
$ws.Range("Q2:Q3").Select()
